$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.102351
$ws.Range("H2").Value = 0.307053
$ws.Range("I2").Value = 0.2080046986044413
$ws.Range("J2").Value = 0.2080046986044413
$ws.Range("M2").Value = 6.382924
$ws.Range("N2").Value = 19.148772
$ws.Range("O2").Value = 0.1363153751023214
$ws.Range("P2").Value = 0.1363153751023214
$ws.Range("Q2").Value = 0.653298654324
$ws.Range("R2").Value = 5.879687888916001
$ws.Range("S2").Value = 0.02835423851330973
$ws.Range("T2").Value = 0.02835423851330973

$ws.Range("G3").Value = 0.102351
$ws.Range("H3").Value = 0.307053
$ws.Range("I3").Value = 0.2080046986044413
$ws.Range("J3").Value = 0.2080046986044413
$ws.Range("O3").Value = 0.6265841681043937
$ws.Range("P3").Value = 0.6265841681043938
$ws.Range("Q3").Value = 3.002937809004
$ws.Range("R3").Value = 27.026440281036
$ws.Range("S3").Value = 0.130332451036869
$ws.Range("T3").Value = 0.130332451036869

$ws.Range("G4").Value = 0.102351
$ws.Range("H4").Value = 0.307053
$ws.Range("I4").Value = 0.2080046986044413
$ws.Range("J4").Value = 0.2080046986044413
$ws.Range("O4").Value = 0.2371004567932849
$ws.Range("P4").Value = 0.2371004567932849
$ws.Range("Q4").Value = 1.13631649582
$ws.Range("R4").Value = 10.22684846238
$ws.Range("S4").Value = 0.04931800905426259
$ws.Range("T4").Value = 0.04931800905426258

$ws.Range("I5").Value = 0.7162147240552154
$ws.Range("J5").Value = 0.7162147240552154
$ws.Range("M5").Value = 6.382924
$ws.Range("N5").Value = 19.148772
$ws.Range("O5").Value = 0.1363153751023214
$ws.Range("P5").Value = 0.1363153751023214
$ws.Range("Q5").Value = 2.249478586645333
$ws.Range("R5").Value = 20.245307279808
$ws.Range("S5").Value = 0.09763107876339229
$ws.Range("T5").Value = 0.0976310787633923

$ws.Range("I6").Value = 0.7162147240552154
$ws.Range("J6").Value = 0.7162147240552154
$ws.Range("O6").Value = 0.6265841681043937
$ws.Range("P6").Value = 0.6265841681043938
$ws.Range("S6").Value = 0.448768807056255
$ws.Range("T6").Value = 0.448768807056255

$ws.Range("I7").Value = 0.7162147240552154
$ws.Range("J7").Value = 0.7162147240552154
$ws.Range("O7").Value = 0.2371004567932849
$ws.Range("P7").Value = 0.2371004567932849
$ws.Range("S7").Value = 0.1698148382355681
$ws.Range("T7").Value = 0.1698148382355681

$ws.Range("G8").Value = 0.03728866666666666
$ws.Range("I8").Value = 0.07578057734034331
$ws.Range("J8").Value = 0.0757805773403433
$ws.Range("M8").Value = 6.382924
$ws.Range("N8").Value = 19.148772
$ws.Range("O8").Value = 0.1363153751023214
$ws.Range("P8").Value = 0.1363153751023214
$ws.Range("Q8").Value = 0.2380107253946667
$ws.Range("R8").Value = 2.142096528552
$ws.Range("S8").Value = 0.01033005782561938
$ws.Range("T8").Value = 0.01033005782561938

$ws.Range("G9").Value = 0.03728866666666666
$ws.Range("I9").Value = 0.07578057734034331
$ws.Range("J9").Value = 0.0757805773403433
$ws.Range("O9").Value = 0.6265841681043937
$ws.Range("P9").Value = 0.6265841681043938
$ws.Range("S9").Value = 0.04748291001126968
$ws.Range("T9").Value = 0.04748291001126968

$ws.Range("G10").Value = 0.03728866666666666
$ws.Range("I10").Value = 0.07578057734034331
$ws.Range("J10").Value = 0.0757805773403433
$ws.Range("O10").Value = 0.2371004567932849
$ws.Range("P10").Value = 0.2371004567932849
$ws.Range("S10").Value = 0.01796760950345425
$ws.Range("T10").Value = 0.01796760950345425

